$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values look like plain numbers (single decimal point).
# Coerce them to Text format first so Excel keeps the exact original string
# (matching digits / trailing zeros) instead of re-parsing them as numbers.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D16", "D17", "D18", "D19", "D21", "D22", "D24", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the refreshed coin values scraped on Mon Apr 17 19:16:00 UTC 2023.
$ws.Range("D2").Value = "29.700.05"
$ws.Range("E2").Value = "  -3.38%  "
$ws.Range("D3").Value = "2.097.88"
$ws.Range("E3").Value = "  -2.36%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "344.05"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").Value = "0.5143"
$ws.Range("E7").Value = "  -2.72%  "
$ws.Range("D8").Value = "0.4410"
$ws.Range("E8").Value = "  -3.38%  "
$ws.Range("D9").Value = "52.75"
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("D10").Value = "0.09186"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").Value = "1.173"
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("D12").Value = "25.01"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").Value = "2.104.16"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").Value = "8.272"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("D16").Value = "99.65"
$ws.Range("E16").Value = "  -2.70%  "
$ws.Range("D17").Value = "0.00001151"
$ws.Range("E17").Value = "  -2.53%  "
$ws.Range("D18").Value = "1.007"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").Value = "20.82"
$ws.Range("E19").Value = "  +6.19%  "
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "6.191"
$ws.Range("E22").Value = "  -3.38%  "
$ws.Range("D23").Value = "29.746.94"
$ws.Range("E23").Value = "  -3.47%  "
$ws.Range("D24").Value = "12.60"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("E25").Value = "  -2.98%  "
$ws.Range("D26").Value = "2.349.32"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").Value = "21.92"
$ws.Range("E27").Value = "  -3.37%  "
$ws.Range("D28").Value = "2.530"
$ws.Range("E28").Value = "  -2.48%  "
$ws.Range("D29").Value = "162.07"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D30").Value = "132.92"
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("D31").Value = "1.131"
$ws.Range("E31").Value = "  -7.11%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.1050"
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "1.660"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("D34").Value = "6.176"
$ws.Range("E34").Value = "  -3.86%  "
$ws.Range("D35").Value = "3.934"
$ws.Range("E35").Value = "  -1.82%  "
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").Value = "6.022"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("D38").Value = "0.02565"
$ws.Range("E38").Value = "  -3.39%  "
$ws.Range("D39").Value = "0.06727"
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").Value = "0.6874"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").Value = "0.2238"
$ws.Range("E42").Value = "  -4.27%  "
$ws.Range("D43").Value = "1.293"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "0.6664"
$ws.Range("E44").Value = "  +2.80%  "
$ws.Range("D45").Value = "14.25"
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("D46").Value = "2.302"
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("E47").Value = "  -4.10%  "
$ws.Range("D48").Value = "0.00000000348"
$ws.Range("E48").Value = "  -5.76%  "
$ws.Range("D49").Value = "1.221"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("D50").Value = "82.29"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").Value = "0.3326"
$ws.Range("E51").Value = "  -2.35%  "
